# Apply latest cryptos snapshot: refresh Price (D) and Volume(1h) (E)
# columns for the rows whose market data moved since the last run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.150.37'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '2.929.17'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'" + '592.34'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = "'" + '144.58'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  +3.82%  '
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = "'" + '0.0000225'
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = "'" + '33.64'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = '3.417.78'
$ws.Range('E15').Value = '  +0.94%  '
$ws.Range('D16').Value = '60.989.00'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = "'" + '6.73'
$ws.Range('E17').Value = '  +0.26%  '
$ws.Range('D18').Value = '2.935.13'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').Value = "'" + '432.80'
$ws.Range('E19').Value = '  +0.84%  '
$ws.Range('D20').Value = "'" + '13.50'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = "'" + '0.680'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('D22').Value = "'" + '7.10'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').Value = "'" + '11.08'
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('D25').Value = "'" + '2.21'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').Value = "'" + '11.86'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').Value = "'" + '0.999'
$ws.Range('D28').Value = "'" + '2.23'
$ws.Range('E28').Value = '  -3.74%  '
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').Value = "'" + '0.111'
$ws.Range('E31').Value = '  +3.56%  '
$ws.Range('D32').Value = "'" + '26.74'
$ws.Range('E32').Value = '  +0.60%  '
$ws.Range('D34').Value = '0.0₃0876'
$ws.Range('E34').Value = '  +3.11%  '
$ws.Range('E35').Value = '  +0.34%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Value = "'" + '2.96'
$ws.Range('E37').Value = '  -2.65%  '
$ws.Range('D38').Value = "'" + '0.123'
$ws.Range('E38').Value = '  -0.25%  '
$ws.Range('D39').Value = "'" + '1.99'
$ws.Range('E39').Value = '  -0.18%  '
$ws.Range('E40').Value = '  +0.43%  '
$ws.Range('D41').Value = "'" + '41.38'
$ws.Range('E41').Value = '  +2.89%  '
$ws.Range('D42').Value = "'" + '0.282'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').Value = "'" + '375.28'
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = "'" + '0.0347'
$ws.Range('E44').Value = '  -0.09%  '
$ws.Range('D45').Value = '2.702.00'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = "'" + '133.27'
$ws.Range('E46').Value = '  +3.05%  '
$ws.Range('D48').Value = "'" + '23.89'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('E51').Value = '  +0.41%  '
